# Auto-generated edit script applying scheduled-runner value refresh
# to the Alpha_Profits workbook (per-sheet Leve profit/price tables).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 24198.584
$ws.Range("I10").Value = 14166.667
$ws.Range("K10").Value = 14166.667
$ws.Range("M10").Value = -13873.667
$ws.Range("H11").Value = 31254.062
$ws.Range("I11").Value = 31254.062
$ws.Range("K11").Value = 31254.062
$ws.Range("M11").Value = -31114.062
$ws.Range("H80").Value = 1017.5
$ws.Range("J80").Value = 1242.1818
$ws.Range("L80").Value = 3726.5454
$ws.Range("N80").Value = -5722.5454
$ws.Range("H83").Value = 1017.5
$ws.Range("J83").Value = 1242.1818
$ws.Range("L83").Value = 11179.6362
$ws.Range("N83").Value = -21163.6362
$ws.Range("H113").Value = 2252.7
$ws.Range("I113").Value = 2065.3125
$ws.Range("K113").Value = 2065.3125
$ws.Range("M113").Value = 1188.6875
$ws.Range("H132").Value = 32294.969
$ws.Range("I132").Value = 35453.31
$ws.Range("K132").Value = 106359.93
$ws.Range("M132").Value = -103829.93

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 1000
$ws.Range("J13").Value = 1000
$ws.Range("L13").Value = 1000
$ws.Range("N13").Value = -1288
$ws.Range("H21").Value = 1099.3334
$ws.Range("I21").Value = 749
$ws.Range("K21").Value = 749
$ws.Range("M21").Value = -375
$ws.Range("H32").Value = 1714.2693
$ws.Range("I32").Value = 1718.2653
$ws.Range("K32").Value = 1718.2653
$ws.Range("M32").Value = -1431.2653
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("H86").Value = 2240.25
$ws.Range("I86").Value = 1514.1666
$ws.Range("J86").Value = 2966.3333
$ws.Range("K86").Value = 1514.1666
$ws.Range("L86").Value = 2966.3333
$ws.Range("M86").Value = -391.1666
$ws.Range("N86").Value = -5212.3333
$ws.Range("H89").Value = 2240.25
$ws.Range("I89").Value = 1514.1666
$ws.Range("J89").Value = 2966.3333
$ws.Range("K89").Value = 7570.833000000001
$ws.Range("L89").Value = 14831.6665
$ws.Range("M89").Value = -1954.833000000001
$ws.Range("N89").Value = -26063.6665
$ws.Range("H134").Value = 15153559
$ws.Range("I134").Value = 2081.4736
$ws.Range("K134").Value = 6244.4208
$ws.Range("M134").Value = -3709.4208

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 20752.75
$ws.Range("I6").Value = 5403.6
$ws.Range("K6").Value = 5403.6
$ws.Range("M6").Value = -5290.6
$ws.Range("H16").Value = 2665.5
$ws.Range("J16").Value = 2832
$ws.Range("L16").Value = 2832
$ws.Range("N16").Value = -3406
$ws.Range("H26").Value = 2919.5454
$ws.Range("I26").Value = 2809.4
$ws.Range("K26").Value = 2809.4
$ws.Range("M26").Value = -2522.4
$ws.Range("H31").Value = 1812.3572
$ws.Range("I31").Value = 1633.9445
$ws.Range("K31").Value = 1633.9445
$ws.Range("M31").Value = -1338.9445
$ws.Range("H34").Value = 1812.3572
$ws.Range("I34").Value = 1633.9445
$ws.Range("K34").Value = 1633.9445
$ws.Range("M34").Value = -1431.9445
$ws.Range("H58").Value = 2221.05
$ws.Range("I58").Value = 2062.2778
$ws.Range("K58").Value = 2062.2778
$ws.Range("M58").Value = -1859.2778
$ws.Range("H80").Value = 17999.666
$ws.Range("J80").Value = 17999.666
$ws.Range("L80").Value = 17999.666
$ws.Range("N80").Value = -20245.666
$ws.Range("H83").Value = 17999.666
$ws.Range("J83").Value = 17999.666
$ws.Range("L83").Value = 53998.99800000001
$ws.Range("N83").Value = -65230.99800000001
$ws.Range("H113").Value = 2665.5
$ws.Range("J113").Value = 2832
$ws.Range("L113").Value = 2832
$ws.Range("N113").Value = -7172
$ws.Range("H134").Value = 20004180
$ws.Range("I134").Value = 3999
$ws.Range("K134").Value = 11997
$ws.Range("M134").Value = -9462
$ws.Range("H136").Value = 2221.05
$ws.Range("I136").Value = 2062.2778
$ws.Range("K136").Value = 6186.8334
$ws.Range("M136").Value = -3636.8334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 371.08694
$ws.Range("I12").Value = 65
$ws.Range("J12").Value = 479.11765
$ws.Range("K12").Value = 195
$ws.Range("L12").Value = 1437.35295
$ws.Range("M12").Value = -22
$ws.Range("N12").Value = -1783.35295
$ws.Range("H113").Value = 552.7692
$ws.Range("I113").Value = 423
$ws.Range("J113").Value = 576.36365
$ws.Range("K113").Value = 1269
$ws.Range("L113").Value = 1729.09095
$ws.Range("M113").Value = 901
$ws.Range("N113").Value = -6069.09095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 6672307.5
$ws.Range("I3").Value = 2500962.5
$ws.Range("K3").Value = 2500962.5
$ws.Range("M3").Value = -2500846.5
$ws.Range("H13").Value = 278.18182
$ws.Range("I13").Value = 128.22223
$ws.Range("K13").Value = 128.22223
$ws.Range("M13").Value = 10.77777
$ws.Range("H22").Value = 5704.5
$ws.Range("J22").Value = 7333.3335
$ws.Range("L22").Value = 7333.3335
$ws.Range("N22").Value = -8391.333500000001
$ws.Range("H102").Value = 2140.5334
$ws.Range("I102").Value = 2401
$ws.Range("K102").Value = 2401
$ws.Range("M102").Value = -779
$ws.Range("H107").Value = 1349.375
$ws.Range("I107").Value = 958.2
$ws.Range("K107").Value = 958.2
$ws.Range("M107").Value = 961.8
$ws.Range("H122").Value = 3306.9443
$ws.Range("I122").Value = 3449.5715
$ws.Range("K122").Value = 10348.7145
$ws.Range("M122").Value = -7898.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 29950
$ws.Range("I4").Value = 29933.334
$ws.Range("K4").Value = 29933.334
$ws.Range("M4").Value = -29820.334
$ws.Range("H12").Value = 14500
$ws.Range("I12").Value = 10000
$ws.Range("K12").Value = 10000
$ws.Range("M12").Value = -9830
$ws.Range("H28").Value = 29950
$ws.Range("I28").Value = 29933.334
$ws.Range("K28").Value = 29933.334
$ws.Range("M28").Value = -29701.334
$ws.Range("H37").Value = 29950
$ws.Range("I37").Value = 29933.334
$ws.Range("K37").Value = 29933.334
$ws.Range("M37").Value = -29826.334
$ws.Range("H62").Value = 45000
$ws.Range("J62").Value = 45000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46248
$ws.Range("H64").Value = 49999.145
$ws.Range("J64").Value = 49999
$ws.Range("L64").Value = 49999
$ws.Range("N64").Value = -50449
$ws.Range("H65").Value = 45000
$ws.Range("J65").Value = 45000
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141240
$ws.Range("H67").Value = 49999.145
$ws.Range("J67").Value = 49999
$ws.Range("L67").Value = 49999
$ws.Range("N67").Value = -51559
$ws.Range("H93").Value = 802.3125
$ws.Range("I93").Value = 664.46155
$ws.Range("K93").Value = 664.46155
$ws.Range("M93").Value = 583.53845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 17614
$ws.Range("I9").Value = 17614
$ws.Range("K9").Value = 17614
$ws.Range("M9").Value = -17474
$ws.Range("H63").Value = 41851.332
$ws.Range("J63").Value = 42777
$ws.Range("L63").Value = 42777
$ws.Range("N63").Value = -44025
$ws.Range("H66").Value = 41851.332
$ws.Range("J66").Value = 42777
$ws.Range("L66").Value = 128331
$ws.Range("N66").Value = -134571
$ws.Range("H96").Value = 3267.348
$ws.Range("J96").Value = 2376.0667
$ws.Range("L96").Value = 2376.0667
$ws.Range("N96").Value = -5122.066699999999

# Explicit cell removals (cells dropped entirely from the row, not just zeroed)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M42").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M10").ClearContents()

Write-Output "Applied updates across sheets: ALC,ARM,BSM,CRP,CUL,GSM,LTW,WVR"